$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Update the date line
$d.Content.Find.Execute("2024-01-19 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-20 Saturday", 2) | Out-Null

# Delete original rows 4-13 (1-indexed), i.e. the 10 rows that get fully replaced
for ($i = 0; $i -lt 10; $i++) {
    $t.Rows.Item(4).Delete()
}

# Set the text for the remaining 10 rows (rows 1-3 kept, rows 4-10 were old rows 14-20)
$r = $t.Rows.Item(1)
$r.Cells.Item(1).Range.Text = "25+33="
$r.Cells.Item(2).Range.Text = "48-37="
$r.Cells.Item(3).Range.Text = "36+0="
$r.Cells.Item(4).Range.Text = "99-84="
$r.Cells.Item(5).Range.Text = "12+49="

$r = $t.Rows.Item(2)
$r.Cells.Item(1).Range.Text = "45-12="
$r.Cells.Item(2).Range.Text = "29+36="
$r.Cells.Item(3).Range.Text = "92-31="
$r.Cells.Item(4).Range.Text = "89-60="
$r.Cells.Item(5).Range.Text = "58+25="

$r = $t.Rows.Item(3)
$r.Cells.Item(1).Range.Text = "70-9="
$r.Cells.Item(2).Range.Text = "80-55="
$r.Cells.Item(3).Range.Text = "12+76="
$r.Cells.Item(4).Range.Text = "64+29="
$r.Cells.Item(5).Range.Text = "24+10="

$r = $t.Rows.Item(4)
$r.Cells.Item(1).Range.Text = "6+49="
$r.Cells.Item(2).Range.Text = "42-24="
$r.Cells.Item(3).Range.Text = "78-78="
$r.Cells.Item(4).Range.Text = "87-39="
$r.Cells.Item(5).Range.Text = "71+14="

$r = $t.Rows.Item(5)
$r.Cells.Item(1).Range.Text = "81-35="
$r.Cells.Item(2).Range.Text = "24+73="
$r.Cells.Item(3).Range.Text = "28+19="
$r.Cells.Item(4).Range.Text = "12+72="
$r.Cells.Item(5).Range.Text = "19+42="

$r = $t.Rows.Item(6)
$r.Cells.Item(1).Range.Text = "54+3="
$r.Cells.Item(2).Range.Text = "31+30="
$r.Cells.Item(3).Range.Text = "53+39="
$r.Cells.Item(4).Range.Text = "70-16="
$r.Cells.Item(5).Range.Text = "73-34="

$r = $t.Rows.Item(7)
$r.Cells.Item(1).Range.Text = "18+39="
$r.Cells.Item(2).Range.Text = "3+46="
$r.Cells.Item(3).Range.Text = "20+60="
$r.Cells.Item(4).Range.Text = "93-32="
$r.Cells.Item(5).Range.Text = "37+29="

$r = $t.Rows.Item(8)
$r.Cells.Item(1).Range.Text = "64-13="
$r.Cells.Item(2).Range.Text = "87-57="
$r.Cells.Item(3).Range.Text = "48+50="
$r.Cells.Item(4).Range.Text = "40+8="
$r.Cells.Item(5).Range.Text = "12+73="

$r = $t.Rows.Item(9)
$r.Cells.Item(1).Range.Text = "62-5="
$r.Cells.Item(2).Range.Text = "88-23="
$r.Cells.Item(3).Range.Text = "47+17="
$r.Cells.Item(4).Range.Text = "13+37="
$r.Cells.Item(5).Range.Text = "4+46="

$r = $t.Rows.Item(10)
$r.Cells.Item(1).Range.Text = "0+1="
$r.Cells.Item(2).Range.Text = "5+71="
$r.Cells.Item(3).Range.Text = "96-7="
$r.Cells.Item(4).Range.Text = "16-3="
$r.Cells.Item(5).Range.Text = "16+57="

# Append 10 brand-new rows at the end
$r = $t.Rows.Add()
$r.Cells.Item(1).Range.Text = "14+77="
$r.Cells.Item(2).Range.Text = "40-19="
$r.Cells.Item(3).Range.Text = "84-9="
$r.Cells.Item(4).Range.Text = "76-12="
$r.Cells.Item(5).Range.Text = "7+40="

$r = $t.Rows.Add()
$r.Cells.Item(1).Range.Text = "11+49="
$r.Cells.Item(2).Range.Text = "56-22="
$r.Cells.Item(3).Range.Text = "35+45="
$r.Cells.Item(4).Range.Text = "92-45="
$r.Cells.Item(5).Range.Text = "2+95="

$r = $t.Rows.Add()
$r.Cells.Item(1).Range.Text = "91-21="
$r.Cells.Item(2).Range.Text = "91-11="
$r.Cells.Item(3).Range.Text = "59-33="
$r.Cells.Item(4).Range.Text = "53+23="
$r.Cells.Item(5).Range.Text = "70-45="

$r = $t.Rows.Add()
$r.Cells.Item(1).Range.Text = "1+88="
$r.Cells.Item(2).Range.Text = "24+39="
$r.Cells.Item(3).Range.Text = "0+54="
$r.Cells.Item(4).Range.Text = "60-7="
$r.Cells.Item(5).Range.Text = "12+20="

$r = $t.Rows.Add()
$r.Cells.Item(1).Range.Text = "47+38="
$r.Cells.Item(2).Range.Text = "8+82="
$r.Cells.Item(3).Range.Text = "73+11="
$r.Cells.Item(4).Range.Text = "0+1="
$r.Cells.Item(5).Range.Text = "10+85="

$r = $t.Rows.Add()
$r.Cells.Item(1).Range.Text = "62-60="
$r.Cells.Item(2).Range.Text = "29+42="
$r.Cells.Item(3).Range.Text = "61-56="
$r.Cells.Item(4).Range.Text = "91-78="
$r.Cells.Item(5).Range.Text = "86-66="

$r = $t.Rows.Add()
$r.Cells.Item(1).Range.Text = "16+52="
$r.Cells.Item(2).Range.Text = "12+82="
$r.Cells.Item(3).Range.Text = "70-21="
$r.Cells.Item(4).Range.Text = "55+30="
$r.Cells.Item(5).Range.Text = "43-0="

$r = $t.Rows.Add()
$r.Cells.Item(1).Range.Text = "29+69="
$r.Cells.Item(2).Range.Text = "0+41="
$r.Cells.Item(3).Range.Text = "24-6="
$r.Cells.Item(4).Range.Text = "85-58="
$r.Cells.Item(5).Range.Text = "23+10="

$r = $t.Rows.Add()
$r.Cells.Item(1).Range.Text = "91-2="
$r.Cells.Item(2).Range.Text = "41+27="
$r.Cells.Item(3).Range.Text = "36+3="
$r.Cells.Item(4).Range.Text = "19-7="
$r.Cells.Item(5).Range.Text = "77-3="

$r = $t.Rows.Add()
$r.Cells.Item(1).Range.Text = "51-0="
$r.Cells.Item(2).Range.Text = "27+4="
$r.Cells.Item(3).Range.Text = "72-59="
$r.Cells.Item(4).Range.Text = "67+20="
$r.Cells.Item(5).Range.Text = "71-29="
